# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
#
# The two oldest pending fixtures (rows 244 and 245, ids 242/243 with
# B=6775592 / B=6774471) are removed from the feed. All later pending
# fixtures shift up two rows, their sequential "id" (column A) is
# renumbered accordingly, and the odds feed refreshed a few of the
# remaining fixtures' quoted prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows; everything below shifts up by two rows.
$ws.Range("A244:A245").EntireRow.Delete()

# Renumber the sequential id column (A) for the shifted rows.
$ws.Cells.Item(244, 1).Value = 242
$ws.Cells.Item(245, 1).Value = 243
$ws.Cells.Item(246, 1).Value = 244
$ws.Cells.Item(247, 1).Value = 245
$ws.Cells.Item(248, 1).Value = 246
$ws.Cells.Item(249, 1).Value = 247
$ws.Cells.Item(250, 1).Value = 248

# Refresh updated odds for the fixtures that changed along with the shift.

# Row 246 (Rakow Czestochowa vs Legia Warsaw)
$ws.Cells.Item(246, 18).Value = 1.825   # R246
$ws.Cells.Item(246, 19).Value = 2.025   # S246

# Row 247 (LKS Lodz vs Radomiak Radom)
$ws.Cells.Item(247, 14).Value = 3.5     # N247
$ws.Cells.Item(247, 16).Value = 2.05    # P247
$ws.Cells.Item(247, 18).Value = 1.825   # R247
$ws.Cells.Item(247, 19).Value = 2.025   # S247

# Row 248 (Jagiellonia Bialystok vs Cracovia Krakow)
$ws.Cells.Item(248, 14).Value = 1.7     # N248
$ws.Cells.Item(248, 16).Value = 4.333   # P248
$ws.Cells.Item(248, 18).Value = 1.925   # R248
$ws.Cells.Item(248, 19).Value = 1.925   # S248
$ws.Cells.Item(248, 20).Value = 2.75    # T248
$ws.Cells.Item(248, 21).Value = 2.025   # U248
$ws.Cells.Item(248, 22).Value = 1.825   # V248

# Row 249 (Gornik Zabrze vs Slask Wroclaw)
$ws.Cells.Item(249, 21).Value = 1.9     # U249
$ws.Cells.Item(249, 22).Value = 1.95    # V249
